$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 80. This shifts the existing rows 80..220
# down to 81..221 (and updates the sheet dimension to A1:R221).
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with the new weekly data point.
$ws.Range('A80').Value = 3
$ws.Range('B80').Value = 'Femacal de La Calera'
$ws.Range('C80').Value = 'Coquimbo'
$ws.Range('D80').Value = 44533
$ws.Range('E80').Value = 5
$ws.Range('F80').Value = 100112039
$ws.Range('G80').Value = 'Ciboulette'
$ws.Range('H80').Value = 'Sin especificar'
$ws.Range('I80').Value = 'Primera'
$ws.Range('J80').Value = 160
$ws.Range('K80').Value = 1500
$ws.Range('L80').Value = 1500
$ws.Range('M80').Value = 1500
$ws.Range('N80').Value = '$/docena de atados'
$ws.Range('O80').Value = 'Provincia de Quillota'
$ws.Range('P80').Value = 500
$ws.Range('Q80').Value = 3
$ws.Range('R80').Value = 'Hortaliza'
